$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tbl = $shape.Table

# Remove the extra "Behavior" row (row 3) - it was an erroneous extra layer of nesting
$tbl.Rows.Item(3).Delete()

# After deletion, rows shift up by one:
# Row 3 is now "Victim_Targeting" (was "    Victim_Targeting")
$cell = $tbl.Cell(3, 1)
$cell.Shape.TextFrame.TextRange.Text = ""

# Row 4 is now "Identity" cell (was "        Identity", 8 leading spaces -> 4 leading spaces)
$cell = $tbl.Cell(4, 1)
$cell.Shape.TextFrame.TextRange.Text = "    "

# Row 5 is now "Specification" cell (was "            Specification", 12 leading spaces -> 8 leading spaces)
$cell = $tbl.Cell(5, 1)
$cell.Shape.TextFrame.TextRange.Text = "        Specification"

# Row 6 is now "OrganisationInfo" cell (was "                OrganisationInfo", 16 leading spaces -> 12 leading spaces)
$cell = $tbl.Cell(6, 1)
$cell.Shape.TextFrame.TextRange.Text = "            "
